# Updated symbol list (coin rankings/prices) to match the
# "Fri Dec 30 17:52:01 UTC 2022" GitHub Actions refresh.
# Numeric-looking values in column D are written with a leading
# apostrophe so Excel keeps them as text (preserving formatting such
# as trailing zeros, e.g. "244.80"), matching the original sheet where
# every Price cell is stored as a string rather than a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'244.80"
$ws.Range('D3').Value = "'24.99"
$ws.Range('B4').Value = 'HuobiToken'
$ws.Range('C4').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D4').Value = "'5.055"
$ws.Range('E4').Value = '3HuobiTokenHT'
$ws.Range('B5').Value = 'Cronos'
$ws.Range('C5').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D5').Value = "'0.05620"
$ws.Range('E5').Value = '4CronosCRO'
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D6').Value = "'6.524"
$ws.Range('E6').Value = '5KuCoinTokenKCS'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = "'2.981"
$ws.Range('E7').Value = '6GateTokenGT'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = "'0.8099"
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = "'0.8398"
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'0.1337"
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = "'0.03332"
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('D12').Value = "'0.06952"
$ws.Range('D13').Value = "'0.02841"
$ws.Range('D14').Value = "'0.09404"
$ws.Range('D15').Value = "'0.001528"
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').Value = "'0.0005985"
$ws.Range('E16').Value = '15OneONE'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = "'0.006232"
$ws.Range('E17').Value = '16TigerCashTCH'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = "'3.498"
$ws.Range('E18').Value = '17LEOLEO'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').Value = "'2.091"
$ws.Range('E19').Value = '18BTSETokenBTSE'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').Value = "'0.3171"
$ws.Range('E20').Value = '19BitpandaEcosystemTokenBEST'
$ws.Range('D22').Value = "'3.757"
$ws.Range('D23').Value = "'0.04678"
$ws.Range('D24').Value = "'0.1369"
$ws.Range('D25').Value = "'0.001243"
$ws.Range('D27').Value = "'0.00009703"
$ws.Range('E27').Value = '26NitroExNTXBestin24h'
$ws.Range('D28').Value = "'0.0001938"
$ws.Range('D40').Value = "'0.03631"
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = "'0.1052"
$ws.Range('E41').Value = '40BKEXTokenBKK'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = "'0.002717"
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = "'0.003366"
$ws.Range('E43').Value = '42KickTokenKICK'
$ws.Range('D44').Value = "'0.008359"
$ws.Range('D45').Value = "'0.00005293"
$ws.Range('D46').Value = "'0.00000000749"
$ws.Range('D47').Value = "'0.1798"
$ws.Range('D48').Value = "'0.002285"
$ws.Range('D49').Value = "'0.00002098"
$ws.Range('D50').Value = "'0.0001998"
